$newValues = @(
    'Sát Phá Tham hội chiếu tại Điền Trạch',
    'Tử Phủ Vũ Tướng hội chiếu tại Điền Trạch',
    'Cơ Nguyệt Đồng Lương hội chiếu tại Điền Trạch',
    'Cự Nhật hội chiếu tại Điền Trạch',
    'Kình Đà hội chiếu tại Điền Trạch',
    'Xương Khúc hội chiếu tại Điền Trạch',
    'Hoả Linh hội chiếu tại Điền Trạch',
    'Không Kiếp hội chiếu tại Điền Trạch',
    'Quang Quý hội chiếu tại Điền Trạch',
    'Tả Hữu hội chiếu tại Điền Trạch',
    'Song Hao hội chiếu tại Điền Trạch',
    'Tang Hổ hội chiếu tại Điền Trạch',
    'Khốc Hư hội chiếu tại Điền Trạch',
    'Hình Riêu hội chiếu tại Điền Trạch',
    'Thai Toạ hội chiếu tại Điền Trạch',
    'Đào Hồng hội chiếu tại Điền Trạch',
    'Ấn Phù hội chiếu tại Điền Trạch',
    'Song Hao Quyền Lộc Kiếp Hoả hội chiếu tại Điền Trạch',
    'Tử Phủ Vũ Tướng Xương Khúc Khôi Việt Tả Hữu Khoa Quyền Lộc Long hội chiếu tại Điền Trạch',
    'Tử Phủ Vũ Tướng Tả Hữu Khoa Quyền Lộc Long Phượng hội chiếu tại Điền Trạch',
    'Tử Khúc Phá Dương Đà hội chiếu tại Điền Trạch',
    'Cơ Nguyệt Đồng Lương gặp Xương Khúc Tả Hữu hội chiếu tại Điền Trạch',
    'Cơ Nguyệt Đồng Lương Khoa Tả Hữu Quang Quý Quan Phúc hội chiếu tại Điền Trạch',
    'Sát Quyền hội chiếu tại Điền Trạch',
    'Lộc Mã hội chiếu tại Điền Trạch',
    'Kiếp Hư Hao Quyền hội chiếu tại Điền Trạch',
    'Tuế Hổ Phù Xương Khúc hội chiếu tại Điền Trạch',
    'Xương Khúc Tấu Long Phượng hội chiếu tại Điền Trạch',
    'Đào Hồng Riêu Tấu Cơ Vũ hội chiếu tại Điền Trạch',
    'Binh Hình Tướng Ấn hội chiếu tại Điền Trạch',
    'Hổ Tấu hội chiếu tại Điền Trạch',
    'Hình Riêu Y hội chiếu tại Điền Trạch',
    'Mã Hỏa Linh hội chiếu tại Điền Trạch',
    'Thai Tọa Hồng Đào hội chiếu tại Điền Trạch',
    'Tả Hữu Không Kiếp hội chiếu tại Điền Trạch',
    'Tả Hữu Binh Tướng hội chiếu tại Điền Trạch',
    'Đào Quyền hội chiếu tại Điền Trạch',
    'Đào Hồng Tả Cái Hữu hội chiếu tại Điền Trạch',
    'Quan Phúc Quang Tấu hội chiếu tại Điền Trạch',
    'Đào Hồng Xương Khúc hội chiếu tại Điền Trạch',
    'Đào Tử Phủ hội chiếu tại Điền Trạch',
    'Hổ Kình Sát hội chiếu tại Điền Trạch',
    'Hổ Tang Kiếp hội chiếu tại Điền Trạch',
    'Hổ Tang Không Kiếp hội chiếu tại Điền Trạch',
    'Âm Dương Lương hội chiếu tại Điền Trạch',
    'Cơ Lương Gia Hội hội chiếu tại Điền Trạch',
    'Nhật Chiếu Lôi Môn hội chiếu tại Điền Trạch',
    'Tả Hữu Xương Khúc hội chiếu tại Điền Trạch',
    'Tham Linh Triều Viên hội chiếu tại Điền Trạch',
    'Cự Hỏa Linh hội chiếu tại Điền Trạch',
    'Hồng Đào Kỵ hội chiếu tại Điền Trạch',
    'Cự Đồng Hình hội chiếu tại Điền Trạch',
    'Lương Phá hội chiếu tại Điền Trạch',
    'Khôi Việt hội chiếu tại Điền Trạch',
    'Tham Vũ Hỏa hội chiếu tại Điền Trạch',
    'Cự Kỵ hội chiếu tại Điền Trạch',
    'Kình Đà Hỏa hội chiếu tại Điền Trạch',
    'Kình Đà Không Kiếp hội chiếu tại Điền Trạch',
    'Đào Hồng Hỷ hội chiếu tại Điền Trạch',
    'Đào Hồng Kiếp Sát hội chiếu tại Điền Trạch',
    'Đào Hồng Kỵ hội chiếu tại Điền Trạch',
    'Đào Hồng Tả Phù Hữu Bật hội chiếu tại Điền Trạch',
    'Đào Hồng Tả Phù Hữu Bật Khoa Quyền Lộc hội chiếu tại Điền Trạch',
    'Đào Hồng Tả Phù Hữu Bật Khoa Quyền Lộc Long Trì Phượng Các hội chiếu tại Điền Trạch',
    'Lưu Hà Kiếp Sát hội chiếu tại Điền Trạch',
    'Phục Không Kiếp hội chiếu tại Điền Trạch',
    'Xương Khúc Khôi Việt hội chiếu tại Điền Trạch',
    'Khoa Quyền hội chiếu tại Điền Trạch',
    'Lộc Quyền hội chiếu tại Điền Trạch',
    'Tướng Binh Đào Hồng hội chiếu tại Điền Trạch',
    'Lương Khốc Tuế hội chiếu tại Điền Trạch',
    'Tướng Binh hội chiếu tại Điền Trạch',
    'Cự Tang hội chiếu tại Điền Trạch',
    'Cự Hỏa hội chiếu tại Điền Trạch',
    'Phá Hình Kỵ hội chiếu tại Điền Trạch',
    'Tang Trực Tuế hội chiếu tại Điền Trạch',
    'Tuế Xương Khúc hội chiếu tại Điền Trạch',
    'Tử Vi Tả Hữu hội chiếu tại Điền Trạch',
    'Tử Vi Tang Tả Hữu hội chiếu tại Điền Trạch'
)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 4247
$count = $newValues.Count

for ($i = 0; $i -lt $count; $i++) {
    $row = $startRow + $i
    $val = $newValues[$i]
    $ws.Cells.Item($row, 1).Value = $val
    $ws.Cells.Item($row, 2).Value = $val
}

$lastRow = $startRow + $count - 1
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4218
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B" + $startRow + ":B" + $lastRow).Select()
